{"js": "// Replace each two-digit multiplication problem in the document with its\n// new value. Every \"old\" string below is unique within the document, so a\n// simple exact-text search-and-replace on the body is sufficient.\nconst replacements = [\n  [\"98\u00d765=\", \"32\u00d788=\"],\n  [\"54\u00d725=\", \"52\u00d735=\"],\n  [\"38\u00d786=\", \"97\u00d776=\"],\n  [\"83\u00d722=\", \"31\u00d755=\"],\n  [\"42\u00d716=\", \"52\u00d778=\"],\n  [\"48\u00d759=\", \"79\u00d795=\"],\n  [\"36\u00d736=\", \"33\u00d751=\"],\n  [\"11\u00d711=\", \"53\u00d752=\"],\n  [\"66\u00d721=\", \"91\u00d747=\"],\n  [\"40\u00d771=\", \"46\u00d799=\"],\n  [\"77\u00d727=\", \"75\u00d786=\"],\n  [\"20\u00d733=\", \"95\u00d733=\"],\n  [\"28\u00d716=\", \"68\u00d782=\"],\n  [\"63\u00d714=\", \"92\u00d766=\"],\n  [\"83\u00d788=\", \"26\u00d737=\"],\n  [\"86\u00d730=\", \"79\u00d755=\"],\n  [\"40\u00d712=\", \"92\u00d732=\"],\n  [\"49\u00d799=\", \"22\u00d795=\"],\n  [\"83\u00d746=\", \"40\u00d783=\"],\n  [\"80\u00d772=\", \"67\u00d745=\"],\n  [\"15\u00d751=\", \"35\u00d743=\"],\n  [\"45\u00d796=\", \"42\u00d717=\"],\n  [\"80\u00d791=\", \"26\u00d748=\"],\n  [\"77\u00d781=\", \"47\u00d771=\"],\n  [\"98\u00d796=\", \"83\u00d773=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each two-digit multiplication problem in the document with its\n# new value. Every \"old\" string is unique within the document, so a single\n# Find/Replace (ReplaceAll) pass per pair is sufficient.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"98\u00d765=\", \"32\u00d788=\"),\n    @(\"54\u00d725=\", \"52\u00d735=\"),\n    @(\"38\u00d786=\", \"97\u00d776=\"),\n    @(\"83\u00d722=\", \"31\u00d755=\"),\n    @(\"42\u00d716=\", \"52\u00d778=\"),\n    @(\"48\u00d759=\", \"79\u00d795=\"),\n    @(\"36\u00d736=\", \"33\u00d751=\"),\n    @(\"11\u00d711=\", \"53\u00d752=\"),\n    @(\"66\u00d721=\", \"91\u00d747=\"),\n    @(\"40\u00d771=\", \"46\u00d799=\"),\n    @(\"77\u00d727=\", \"75\u00d786=\"),\n    @(\"20\u00d733=\", \"95\u00d733=\"),\n    @(\"28\u00d716=\", \"68\u00d782=\"),\n    @(\"63\u00d714=\", \"92\u00d766=\"),\n    @(\"83\u00d788=\", \"26\u00d737=\"),\n    @(\"86\u00d730=\", \"79\u00d755=\"),\n    @(\"40\u00d712=\", \"92\u00d732=\"),\n    @(\"49\u00d799=\", \"22\u00d795=\"),\n    @(\"83\u00d746=\", \"40\u00d783=\"),\n    @(\"80\u00d772=\", \"67\u00d745=\"),\n    @(\"15\u00d751=\", \"35\u00d743=\"),\n    @(\"45\u00d796=\", \"42\u00d717=\"),\n    @(\"80\u00d791=\", \"26\u00d748=\"),\n    @(\"77\u00d781=\", \"47\u00d771=\"),\n    @(\"98\u00d796=\", \"83\u00d773=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
